$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 235 (TAS / Tashkent, Uzbekistan) entirely; all rows below shift up by one.
$ws.Rows.Item(235).Delete()
